# Update cryptos list with latest prices / volume(1h) percentages
# (commit: "Updated cryptos list on Fri May 24 16:52:43 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text values in this sheet
# (e.g. "3.672.41" is not a valid number, and percentages carry padding
# spaces), so force text format before writing to stop COM from silently
# coercing plain-numeric-looking strings (e.g. "596.33") into numbers;
# ClearFormats afterwards restores the original (unstyled) look of the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '68.260.98'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '3.662.16'
$ws.Range('E3').Value = '  -3.28%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '596.33'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '165.66'
$ws.Range('E6').Value = '  -4.10%  '
$ws.Range('D7').Value = '3.660.74'
$ws.Range('E7').Value = '  -3.27%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +3.07%  '
$ws.Range('D11').Value = '6.27'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '37.73'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '4.281.68'
$ws.Range('E15').Value = '  -3.53%  '
$ws.Range('D16').Value = '3.658.64'
$ws.Range('E16').Value = '  -3.73%  '
$ws.Range('D17').Value = '68.116.68'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '7.20'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').Value = '17.01'
$ws.Range('E20').Value = '  +6.19%  '
$ws.Range('D21').Value = '489.62'
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').Value = '9.05'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').Value = '0.717'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '84.24'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  -4.25%  '
$ws.Range('D27').Value = '12.13'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('D31').Value = '7.83'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = '2.37'
$ws.Range('E32').Value = '  -2.54%  '
$ws.Range('D33').Value = '31.10'
$ws.Range('E33').Value = '  -4.65%  '
$ws.Range('D34').Value = '3.807.17'
$ws.Range('E34').Value = '  -3.20%  '
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').Value = '3.609.33'
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').Value = '0.989'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').Value = '5.71'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = '0.132'
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('D41').Value = '0.319'
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '48.76'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '427.53'
$ws.Range('E43').Value = '  -5.66%  '
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D48').Value = '40.13'
$ws.Range('E48').Value = '  -3.04%  '
$ws.Range('D49').Value = '141.38'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '0.0347'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('D51').Value = '2.716.11'
$ws.Range('E51').Value = '  -3.81%  '

$dataRange.ClearFormats()
